$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -21.29090000000003
$ws.Range("A14").Value = -20.61249999999998
$ws.Range("A16").Value = -20.31759999999999
$ws.Range("A21").Value = -21.1937
$ws.Range("A23").Value = -21.32760000000003
$ws.Range("A25").Value = -22.37960000000003
